$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 9.930652333333333
$ws.Range("H2").Value = 29.791957
$ws.Range("I2").Value = 0.9673539331442913
$ws.Range("J2").Value = 0.9673539331442912
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 117.044563
$ws.Range("N2").Value = 351.133689
$ws.Range("O2").Value = 0.3245365645427815
$ws.Range("P2").Value = 0.3245365645427815
$ws.Range("Q2").Value = 1162.32886265993
$ws.Range("R2").Value = 10460.95976393937
$ws.Range("S2").Value = 0.3139417221595958
$ws.Range("T2").Value = 0.3139417221595958
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 9.930652333333333
$ws.Range("H3").Value = 29.791957
$ws.Range("I3").Value = 0.9673539331442913
$ws.Range("J3").Value = 0.9673539331442912
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 101.5800373333333
$ws.Range("N3").Value = 304.740112
$ws.Range("O3").Value = 0.281657135515876
$ws.Range("P3").Value = 0.281657135515876
$ws.Range("Q3").Value = 1008.756034764354
$ws.Range("R3").Value = 9078.804312879185
$ws.Range("S3").Value = 0.2724621378394373
$ws.Range("T3").Value = 0.2724621378394372
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 9.930652333333333
$ws.Range("H4").Value = 29.791957
$ws.Range("I4").Value = 0.9673539331442913
$ws.Range("J4").Value = 0.9673539331442912
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 142.0267893333333
$ws.Range("N4").Value = 426.080368
$ws.Range("O4").Value = 0.3938062999413425
$ws.Range("P4").Value = 0.3938062999413425
$ws.Range("Q4").Value = 1410.418666888908
$ws.Range("R4").Value = 12693.76800200018
$ws.Range("S4").Value = 0.3809500731452582
$ws.Range("T4").Value = 0.3809500731452581
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.3351376666666666
$ws.Range("H5").Value = 1.005413
$ws.Range("I5").Value = 0.03264606685570878
$ws.Range("J5").Value = 0.03264606685570878
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 117.044563
$ws.Range("N5").Value = 351.133689
$ws.Range("O5").Value = 0.3245365645427815
$ws.Range("P5").Value = 0.3245365645427815
$ws.Range("Q5").Value = 39.22604173983966
$ws.Range("R5").Value = 353.0343756585569
$ws.Range("S5").Value = 0.01059484238318569
$ws.Range("T5").Value = 0.01059484238318569
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.3351376666666666
$ws.Range("H6").Value = 1.005413
$ws.Range("I6").Value = 0.03264606685570878
$ws.Range("J6").Value = 0.03264606685570878
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 101.5800373333333
$ws.Range("N6").Value = 304.740112
$ws.Range("O6").Value = 0.281657135515876
$ws.Range("P6").Value = 0.281657135515876
$ws.Range("Q6").Value = 34.04329669180622
$ws.Range("R6").Value = 306.389670226256
$ws.Range("S6").Value = 0.009194997676438717
$ws.Range("T6").Value = 0.009194997676438714
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.3351376666666666
$ws.Range("H7").Value = 1.005413
$ws.Range("I7").Value = 0.03264606685570878
$ws.Range("J7").Value = 0.03264606685570878
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 142.0267893333333
$ws.Range("N7").Value = 426.080368
$ws.Range("O7").Value = 0.3938062999413425
$ws.Range("P7").Value = 0.3938062999413425
$ws.Range("Q7").Value = 47.59852678133155
$ws.Range("R7").Value = 428.386741031984
$ws.Range("S7").Value = 0.01285622679608437
$ws.Range("T7").Value = 0.01285622679608437
